$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the indicator text in B4 — the trailing period was replaced by a
# trailing space when the author retyped it.
$cell = $ws.Range("B4")
$cell.Value = "2.а.2: Total official flows (official development assistance plus other official flows) to the agriculture sector "

# Re-apply the font so the cell picks up its own font record, matching the
# author's re-entry of the text (creates a dedicated style for this cell).
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

# Move the active selection to B7, matching the saved view state.
$ws.Range("B7").Select()
